# Apply cryptos list price/volume update (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.794.90"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.247.07"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'317.62"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'101.34"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "'0.577"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.556"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "'36.94"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Value = "'0.0831"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "'7.68"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "2.590.83"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "'0.856"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "'14.20"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").Value = "2.250.23"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "43.737.13"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "'13.37"
$ws.Range("E19").Value = "  -8.02%  "
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'6.54"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'65.60"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'3.13"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "'235.29"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'2.12"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'10.10"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'37.10"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  -4.95%  "
$ws.Range("D30").Value = "'6.23"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").Value = "'158.83"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").Value = "'20.14"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").Value = "'0.0849"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").Value = "'0.116"
$ws.Range("E35").Value = "  +11.08%  "
$ws.Range("D36").Value = "'3.08"
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").Value = "'4.23"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("D41").Value = "'15.85"
$ws.Range("E41").Value = "  +19.12%  "
$ws.Range("D42").Value = "'0.0316"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "1.789.45"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D45").Value = "'0.198"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").Value = "'82.72"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("D47").Value = "'75.51"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'5.19"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").Value = "'58.52"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'103.54"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.67"
$ws.Range("E51").Value = "  +4.45%  "
